# Auto-generated edit script: updates cryptos price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.787.10'
$ws.Range('E2').Value = '  -0.21%  '
$ws.Range('D3').Value = '1.637.52'
$ws.Range('E3').Value = '  -0.18%  '
$ws.Range('E4').Value = '  -0.21%  '
$ws.Range('D5').Value = "'216.72"
$ws.Range('E5').Value = '  -1.16%  '
$ws.Range('E6').Value = '  +1.87%  '
$ws.Range('E7').Value = '  -0.26%  '
$ws.Range('E8').Value = '  +1.49%  '
$ws.Range('E9').Value = '  +0.07%  '
$ws.Range('D10').Value = "'19.86"
$ws.Range('E10').Value = '  +3.30%  '
$ws.Range('E11').Value = '  +0.13%  '
$ws.Range('D12').Value = '1.867.06'
$ws.Range('E12').Value = '  +0.06%  '
$ws.Range('D13').Value = '1.638.55'
$ws.Range('E13').Value = '  -0.33%  '
$ws.Range('E14').Value = '  -0.52%  '
$ws.Range('E15').Value = '  +0.78%  '
$ws.Range('D16').Value = "'66.39"
$ws.Range('E16').Value = '  +2.71%  '
$ws.Range('D17').Value = '26.792.68'
$ws.Range('E17').Value = '  -0.01%  '
$ws.Range('E18').Value = '  -0.74%  '
$ws.Range('D19').Value = "'218.86"
$ws.Range('E19').Value = '  +1.50%  '
$ws.Range('E20').Value = '  -0.44%  '
$ws.Range('D21').Value = "'6.68"
$ws.Range('E21').Value = '  +6.50%  '
$ws.Range('E22').Value = '  +0.77%  '
$ws.Range('E23').Value = '  +3.94%  '
$ws.Range('E24').Value = '  +0.38%  '
$ws.Range('D25').Value = "'147.16"
$ws.Range('E25').Value = '  -0.38%  '
$ws.Range('E26').Value = '  -0.22%  '
$ws.Range('E27').Value = '  +4.38%  '
$ws.Range('E28').Value = '  +0.41%  '
$ws.Range('D29').Value = "'15.73"
$ws.Range('E29').Value = '  +0.22%  '
$ws.Range('E30').Value = '  -0.47%  '
$ws.Range('E31').Value = '  -1.87%  '
$ws.Range('E32').Value = '  -1.95%  '
$ws.Range('D33').Value = "'3.00"
$ws.Range('E33').Value = '  +0.60%  '
$ws.Range('E34').Value = '  +1.25%  '
$ws.Range('D35').Value = '1.259.50'
$ws.Range('E35').Value = '  -0.17%  '
$ws.Range('E36').Value = '  +0.08%  '
$ws.Range('E37').Value = '  +1.15%  '
$ws.Range('D38').Value = "'0.533"
$ws.Range('E38').Value = '  +0.98%  '
$ws.Range('D39').Value = "'0.833"
$ws.Range('E39').Value = '  +2.63%  '
$ws.Range('E40').Value = '  -0.44%  '
$ws.Range('D41').Value = "'0.806"
$ws.Range('E41').Value = '  -0.04%  '
$ws.Range('E42').Value = '  +2.59%  '
$ws.Range('D43').Value = '1.777.31'
$ws.Range('E43').Value = '  +0.20%  '
$ws.Range('B44').Value = 'MXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D44').Value = "'2.11"
$ws.Range('E44').Value = '  -0.78%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').Value = "'61.72"
$ws.Range('E45').Value = '  +2.87%  '
$ws.Range('D46').Value = "'91.67"
$ws.Range('E46').Value = '  -0.49%  '
$ws.Range('E47').Value = '  -1.30%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').Value = "'0.0513"
$ws.Range('E48').Value = '  -0.84%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').Value = "'7.63"
$ws.Range('E49').Value = '  +1.58%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').Value = "'0.0964"
$ws.Range('E50').Value = '  +0.34%  '
$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D51').Value = "'1.00"
$ws.Range('E51').Value = '  -0.30%  '
